$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the shared color strings (B2 and C2 swap which shared-string index they
# reference, and the three color strings themselves are renamed).
$ws.Range("B2").Value = "#4285F4"
$ws.Range("C2").Value = "#DB4437"
$ws.Range("D2").Value = "#F4B400"

# Update selection from A3 to E2
$ws.Range("E2").Select()
